# ST1 Tool Monitoring log update
#  1) Row 2: store the tool-life "cnt" / sign-off readings as real numbers (were text)
#  2) Row 3: append the next shift's monitoring entry
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: convert the numeric reading columns from text to numbers ---
$ws.Range("D2").Value = 1
$ws.Range("F2").Value = 341
$ws.Range("H2").Value = 31
$ws.Range("J2").Value = 321
$ws.Range("L2").Value = 321
$ws.Range("N2").Value = 321
$ws.Range("P2").Value = 321
$ws.Range("R2").Value = 32
$ws.Range("T2").Value = 3
$ws.Range("V2").Value = 315
$ws.Range("X2").Value = 31
$ws.Range("Z2").Value = 312
$ws.Range("AB2").Value = 321
$ws.Range("AD2").Value = 123
$ws.Range("AF2").Value = 123
$ws.Range("AG2").Value = 231
$ws.Range("AH2").Value = 123
$ws.Range("AI2").Value = 132

# --- Row 3: new monitoring entry ---
# The digit-only readings must be kept as text (same convention as the rest of the log),
# so format those cells as Text first - otherwise Excel auto-converts digit strings to numbers.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("N3").NumberFormat = "@"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("R3").NumberFormat = "@"
$ws.Range("T3").NumberFormat = "@"
$ws.Range("V3").NumberFormat = "@"
$ws.Range("X3").NumberFormat = "@"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AD3").NumberFormat = "@"
$ws.Range("AF3").NumberFormat = "@"
$ws.Range("AG3").NumberFormat = "@"
$ws.Range("AH3").NumberFormat = "@"
$ws.Range("AI3").NumberFormat = "@"

$ws.Range("A3").Value = "2025-02-04T17:28"
$ws.Range("B3").Value = "SHIFT1"
$ws.Range("C3").Value = "OK"
$ws.Range("D3").Value = "4"
$ws.Range("E3").Value = "OK"
$ws.Range("F3").Value = "4"
$ws.Range("G3").Value = "OK"
$ws.Range("H3").Value = "59"
$ws.Range("I3").Value = "OK"
$ws.Range("J3").Value = "4"
$ws.Range("K3").Value = "OK"
$ws.Range("L3").Value = "41"
$ws.Range("M3").Value = "OK"
$ws.Range("N3").Value = "4"
$ws.Range("O3").Value = "OK"
$ws.Range("P3").Value = "4"
$ws.Range("Q3").Value = "OK"
$ws.Range("R3").Value = "4"
$ws.Range("S3").Value = "OK"
$ws.Range("T3").Value = "4"
$ws.Range("U3").Value = "OK"
$ws.Range("V3").Value = "4"
$ws.Range("W3").Value = "OK"
$ws.Range("X3").Value = "4"
$ws.Range("Y3").Value = "OK"
$ws.Range("Z3").Value = "4"
$ws.Range("AA3").Value = "OK"
$ws.Range("AB3").Value = "4"
$ws.Range("AC3").Value = "OK"
$ws.Range("AD3").Value = "5"
$ws.Range("AE3").Value = "OK"
$ws.Range("AF3").Value = "4"
$ws.Range("AG3").Value = "4"
$ws.Range("AH3").Value = "4"
$ws.Range("AI3").Value = "4"
